$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "<down>"
$ws.Range("B2").Value = "<the>"
$ws.Range("C2").Value = 28

$ws.Range("A3").Value = "<is>"
$ws.Range("B3").Value = "<is>"
$ws.Range("C3").Value = 32

$ws.Range("A4").Value = "<not>"
$ws.Range("B4").Value = "<not>"
$ws.Range("C4").Value = 38

$ws.Range("A5").Value = "<three>"
$ws.Range("B5").Value = "<three>"
$ws.Range("C5").Value = 33

$ws.Range("A6").Value = "<it>"
$ws.Range("B6").Value = "<it>"
$ws.Range("C6").Value = 33

$ws.Range("A7").Value = "<each>"
$ws.Range("B7").Value = "<each>"
$ws.Range("C7").Value = 32

$ws.Range("A8").Value = "<lima>"
$ws.Range("B8").Value = "<number>"
$ws.Range("C8").Value = 36

$ws.Range("A9").Value = "<foxtrot>"
$ws.Range("B9").Value = "<come>"
$ws.Range("C9").Value = 32

$ws.Range("A10").Value = "<a>"
$ws.Range("B10").Value = "<day>"
$ws.Range("C10").Value = 29

$ws.Range("A11").Value = "<and>"
$ws.Range("B11").Value = "<cape>"
$ws.Range("C11").Value = 37

$ws.Range("A12").Value = "<is>"
$ws.Range("B12").Value = "<is>"
$ws.Range("C12").Value = 29

$ws.Range("A13").Value = "<five>"
$ws.Range("B13").Value = "<five>"
$ws.Range("C13").Value = 35

$ws.Range("A14").Value = "<november>"
$ws.Range("B14").Value = "<november>"
$ws.Range("C14").Value = 36

$ws.Range("A15").Value = "<nine>"
$ws.Range("B15").Value = "<nine>"
$ws.Range("C15").Value = 32

$ws.Range("A16").Value = "<number>"
$ws.Range("B16").Value = "<number>"
$ws.Range("C16").Value = 30

$ws.Range("A17").Value = "<escape>"
$ws.Range("B17").Value = "<entee>"
$ws.Range("C17").Value = 35

$ws.Range("A18").Value = "<but>"
$ws.Range("B18").Value = "<not>"
$ws.Range("C18").Value = 23
